$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 628; this shifts existing rows 628..711 down to 629..712
$ws.Rows.Item(628).Insert()

# Populate the new row 628 with the new weekly record
$ws.Cells.Item(628, 1).Value = 10
$ws.Cells.Item(628, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(628, 3).Value = "La Araucanía"
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(628, 4).Value = $epoch.AddDays(44776)
$ws.Cells.Item(628, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(628, 5).Value = 9
$ws.Cells.Item(628, 6).Value = 100112003
$ws.Cells.Item(628, 7).Value = "Ajo"
$ws.Cells.Item(628, 8).Value = "Chino"
$ws.Cells.Item(628, 9).Value = "Primera"
$ws.Cells.Item(628, 10).Value = 300
$ws.Cells.Item(628, 11).Value = 30000
$ws.Cells.Item(628, 12).Value = 30000
$ws.Cells.Item(628, 13).Value = 30000
$ws.Cells.Item(628, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(628, 15).Value = "China"
$ws.Cells.Item(628, 16).Value = 3000
$ws.Cells.Item(628, 17).Value = 10
$ws.Cells.Item(628, 18).Value = "Hortaliza"
